$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "http://purl.obolibrary.org/obo/BFO_0000015"
$ws.Range("B2").Value = "process"
$ws.Range("C2").Value = "R:only×3"
$ws.Range("D2").Value = "SubClassOf: ns1:BFO_0000003 | SubClassOf: ns1:BFO_0000117 only (ns1:BFO_0000015 OR ns1:BFO_0000035) | SubClassOf: ns1:BFO_0000132 only ns1:BFO_0000015 | SubClassOf: ns1:BFO_0000139 only ns1:BFO_0000015"
$ws.Range("E2").Value = "http://ies.data.gov.uk/ontology/ies4#Event"
$ws.Range("F2").Value = "Event"
$ws.Range("G2").Value = "R:only"
$ws.Range("H2").Value = "SubClassOf: ies:Element | SubClassOf: ns1:BFO_0000015 | SubClassOf: ns1:BFO_0000178 only (ns1:BFO_0000029 OR ns1:BFO_0000140)"
$ws.Range("I2").Value = "(Elucidation) p is a process means p is an occurrent that has some temporal proper part and for some time t, p has some material entity as participant"
$ws.Range("J2").Value = "An Event represents an activity or incident, involving one or more participating entities, that occurred/started at a specific point in time – e.g. a meeting, or a telephone call."
